$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gained a new "2022" column (T), one cell per existing data row,
# carrying the same look as the preceding "2021" column (S).
$ws.Range("S4:S14").Copy()
$ws.Range("T4:T14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("T4").Value = 2022
$ws.Range("T5").Value = 3.7
$ws.Range("T6").Value = 1.6
$ws.Range("T7").Value = 1.7
$ws.Range("T8").Value = 17.899999999999999
$ws.Range("T9").Value = 7.5
$ws.Range("T10").Value = 1.1000000000000001
$ws.Range("T11").Value = 4.4000000000000004
$ws.Range("T12").Value = 3
$ws.Range("T13").Value = 4.0999999999999996
$ws.Range("T14").Value = 0.8

# The saved view now has the new column's header cell selected.
$ws.Range("U4").Select()
